$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown")

$ws.Range("A3").Value = 41695
$ws.Range("E3").Formula = "=19+21+10+12+24+29"
$ws.Range("F3").Formula = "=45+40+25+30+40+30"
$ws.Range("A4:G7").ClearContents()
$ws.Range("G3").Formula = "=(E3-E2)/F3*60"

$excel.CalculateFullRebuild()

Write-Host ("G3 text = {0}" -f $ws.Range("G3").Text)
